# Corrected IFRS list values (units rebased, data previously mis-scaled)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 148590
$ws.Range("E2").Value = 3509
$ws.Range("F2").Value = 3509
$ws.Range("G2").Value = 2773
$ws.Range("H2").Value = 1437
$ws.Range("I2").Value = 1469
$ws.Range("J2").Value = -32
$ws.Range("K2").Value = 103227
$ws.Range("L2").Value = 38537
$ws.Range("M2").Value = 64690
$ws.Range("N2").Value = 64371
$ws.Range("O2").Value = 319
$ws.Range("P2").Value = 1714
$ws.Range("Q2").Value = 3085
$ws.Range("R2").Value = -3906
$ws.Range("S2").Value = 476
$ws.Range("T2").Value = 2910
$ws.Range("U2").Value = 175
$ws.Range("V2").Value = 26038
$ws.Range("W2").Value = 2.36
$ws.Range("X2").Value = 0.97
$ws.Range("Y2").Value = 2.31
$ws.Range("Z2").Value = 1.37
$ws.Range("AA2").Value = 59.57
$ws.Range("AB2").Value = 3657.02
$ws.Range("AC2").Value = 4285
$ws.Range("AD2").Value = 37.34
$ws.Range("AE2").Value = 191057
$ws.Range("AF2").Value = 0.84
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 0.63
$ws.Range("AI2").Value = 22.94
$ws.Range("AJ2").Value = 34275419

# Row 3
$ws.Range("D3").Value = 117133
$ws.Range("E3").Value = 16111
$ws.Range("F3").Value = 16111
$ws.Range("G3").Value = 14214
$ws.Range("H3").Value = 9907
$ws.Range("I3").Value = 9925
$ws.Range("J3").Value = -19
$ws.Range("K3").Value = 114678
$ws.Range("L3").Value = 39122
$ws.Range("M3").Value = 75556
$ws.Range("N3").Value = 75248
$ws.Range("O3").Value = 308
$ws.Range("P3").Value = 1714
$ws.Range("Q3").Value = 25956
$ws.Range("R3").Value = -12735
$ws.Range("S3").Value = -3339
$ws.Range("T3").Value = 3554
$ws.Range("U3").Value = 22403
$ws.Range("V3").Value = 24123
$ws.Range("W3").Value = 13.76
$ws.Range("X3").Value = 8.460000000000001
$ws.Range("Y3").Value = 14.22
$ws.Range("Z3").Value = 9.09
$ws.Range("AA3").Value = 51.78
$ws.Range("AB3").Value = 4210.52
$ws.Range("AC3").Value = 28957
$ws.Range("AD3").Value = 8.41
$ws.Range("AE3").Value = 223339
$ws.Range("AF3").Value = 1.09
$ws.Range("AG3").Value = 2500
$ws.Range("AH3").Value = 1.03
$ws.Range("AI3").Value = 8.49
$ws.Range("AJ3").Value = 34275419

# Row 4
$ws.Range("D4").Value = 132235
$ws.Range("E4").Value = 25443
$ws.Range("F4").Value = 25443
$ws.Range("G4").Value = 24874
$ws.Range("H4").Value = 18372
$ws.Range("I4").Value = 18358
$ws.Range("J4").Value = 14
$ws.Range("K4").Value = 158668
$ws.Range("L4").Value = 64660
$ws.Range("M4").Value = 94008
$ws.Range("N4").Value = 93631
$ws.Range("O4").Value = 377
$ws.Range("P4").Value = 1714
$ws.Range("Q4").Value = 27006
$ws.Range("R4").Value = -35646
$ws.Range("S4").Value = 11005
$ws.Range("T4").Value = 15881
$ws.Range("U4").Value = 11125
$ws.Range("V4").Value = 41854
$ws.Range("W4").Value = 19.24
$ws.Range("X4").Value = 13.89
$ws.Range("Y4").Value = 21.74
$ws.Range("Z4").Value = 13.44
$ws.Range("AA4").Value = 68.78
$ws.Range("AB4").Value = 5231.49
$ws.Range("AC4").Value = 53561
$ws.Range("AD4").Value = 6.89
$ws.Range("AE4").Value = 277903
$ws.Range("AF4").Value = 1.33
$ws.Range("AG4").Value = 4000
$ws.Range("AH4").Value = 1.08
$ws.Range("AI4").Value = 7.34
$ws.Range("AJ4").Value = 34275419

# Row 5
$ws.Range("D5").Value = 158745
$ws.Range("E5").Value = 29297
$ws.Range("F5").Value = 29297
$ws.Range("G5").Value = 30847
$ws.Range("H5").Value = 22846
$ws.Range("I5").Value = 22439
$ws.Range("J5").Value = 407
$ws.Range("K5").Value = 195510
$ws.Range("L5").Value = 72962
$ws.Range("M5").Value = 122548
$ws.Range("N5").Value = 114896
$ws.Range("O5").Value = 7651
$ws.Range("P5").Value = 1714
$ws.Range("Q5").Value = 31290
$ws.Range("R5").Value = -47176
$ws.Range("S5").Value = 11145
$ws.Range("T5").Value = 20199
$ws.Range("U5").Value = 11091
$ws.Range("V5").Value = 42277
$ws.Range("W5").Value = 18.46
$ws.Range("X5").Value = 14.39
$ws.Range("Y5").Value = 21.52
$ws.Range("Z5").Value = 12.9
$ws.Range("AA5").Value = 59.54
$ws.Range("AB5").Value = 6688.68
$ws.Range("AC5").Value = 65466
$ws.Range("AD5").Value = 5.62
$ws.Range("AE5").Value = 335215
$ws.Range("AF5").Value = 1.1
$ws.Range("AG5").Value = 10500
$ws.Range("AH5").Value = 2.85
$ws.Range("AI5").Value = 16.04
$ws.Range("AJ5").Value = 34275419

# Row 6
$ws.Range("D6").Value = 165450
$ws.Range("E6").Value = 19674
$ws.Range("F6").Value = 19674
$ws.Range("G6").Value = 22527
$ws.Range("H6").Value = 16419
$ws.Range("I6").Value = 15792
$ws.Range("K6").Value = 207991
$ws.Range("L6").Value = 72548
$ws.Range("M6").Value = 135443
$ws.Range("N6").Value = 127347
$ws.Range("P6").Value = 1714
$ws.Range("Q6").Value = 13809
$ws.Range("R6").Value = -17631
$ws.Range("S6").Value = 448
$ws.Range("T6").Value = 18464
$ws.Range("U6").Value = -4655
$ws.Range("V6").Value = 48297
$ws.Range("W6").Value = 11.89
$ws.Range("X6").Value = 9.92
$ws.Range("Y6").Value = 13.04
$ws.Range("Z6").Value = 8.140000000000001
$ws.Range("AA6").Value = 53.56
$ws.Range("AB6").Value = 7390.44
$ws.Range("AC6").Value = 46074
$ws.Range("AD6").Value = 6.01
$ws.Range("AE6").Value = 371541
$ws.Range("AF6").Value = 0.75
$ws.Range("AG6").Value = 10500
$ws.Range("AH6").Value = 3.79
$ws.Range("AI6").Value = 22.79
$ws.Range("AJ6").Value = 34275419

# Row 7
$ws.Range("D7").Value = 153131
$ws.Range("E7").Value = 11151
$ws.Range("G7").Value = 12217
$ws.Range("H7").Value = 8338
$ws.Range("I7").Value = 8025
$ws.Range("K7").Value = 214858
$ws.Range("L7").Value = 73400
$ws.Range("M7").Value = 141458
$ws.Range("N7").Value = 133146
$ws.Range("P7").Value = 1711
$ws.Range("Q7").Value = 14968
$ws.Range("R7").Value = -2484
$ws.Range("S7").Value = -5546
$ws.Range("T7").Value = 12515
$ws.Range("U7").Value = 4710
$ws.Range("W7").Value = 7.28
$ws.Range("X7").Value = 5.45
$ws.Range("Y7").Value = 6.16
$ws.Range("Z7").Value = 3.94
$ws.Range("AA7").Value = 51.89
$ws.Range("AC7").Value = 23413
$ws.Range("AD7").Value = 8.09
$ws.Range("AE7").Value = 388459
$ws.Range("AF7").Value = 0.49
$ws.Range("AG7").Value = 7275
$ws.Range("AH7").Value = 3.84
$ws.Range("AI7").Value = 31.07

# Row 8
$ws.Range("D8").Value = 152726
$ws.Range("E8").Value = 9858
$ws.Range("G8").Value = 11986
$ws.Range("H8").Value = 8653
$ws.Range("I8").Value = 8285
$ws.Range("K8").Value = 220806
$ws.Range("L8").Value = 73214
$ws.Range("M8").Value = 147592
$ws.Range("N8").Value = 138894
$ws.Range("P8").Value = 1711
$ws.Range("Q8").Value = 15372
$ws.Range("R8").Value = -13042
$ws.Range("S8").Value = -3727
$ws.Range("T8").Value = 13665
$ws.Range("U8").Value = 2037
$ws.Range("W8").Value = 6.45
$ws.Range("X8").Value = 5.67
$ws.Range("Y8").Value = 6.09
$ws.Range("Z8").Value = 3.97
$ws.Range("AA8").Value = 49.61
$ws.Range("AC8").Value = 24173
$ws.Range("AD8").Value = 7.84
$ws.Range("AE8").Value = 405229
$ws.Range("AF8").Value = 0.47
$ws.Range("AG8").Value = 7700
$ws.Range("AH8").Value = 4.06
$ws.Range("AI8").Value = 31.85

# Row 9
$ws.Range("D9").Value = 155306
$ws.Range("E9").Value = 11188
$ws.Range("G9").Value = 13435
$ws.Range("H9").Value = 9765
$ws.Range("I9").Value = 9356
$ws.Range("K9").Value = 228884
$ws.Range("L9").Value = 74457
$ws.Range("M9").Value = 154426
$ws.Range("N9").Value = 145359
$ws.Range("P9").Value = 1711
$ws.Range("Q9").Value = 16372
$ws.Range("R9").Value = -13086
$ws.Range("S9").Value = -2156
$ws.Range("T9").Value = 13858
$ws.Range("U9").Value = 3389
$ws.Range("W9").Value = 7.2
$ws.Range("X9").Value = 6.29
$ws.Range("Y9").Value = 6.58
$ws.Range("Z9").Value = 4.34
$ws.Range("AA9").Value = 48.22
$ws.Range("AC9").Value = 27297
$ws.Range("AD9").Value = 6.94
$ws.Range("AE9").Value = 424090
$ws.Range("AF9").Value = 0.45
$ws.Range("AG9").Value = 8315
$ws.Range("AH9").Value = 4.39
$ws.Range("AI9").Value = 30.46

Write-Output "Updated IFRS rows 2-9"
